# Natmi following Dr Hou advice
# Rebuild the LR-pairs sheet so that every (Sending cluster x Target cluster)
# combination for the Efnb2-Epha3 pair is present (3x3 = 9 data rows instead
# of the previous 3), with refreshed statistics for each combination.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efnb2"
$ws.Cells.Item(2, 3).Value = "Epha3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = [double]"36.899643"
$ws.Cells.Item(2, 8).Value = [double]"110.698929"
$ws.Cells.Item(2, 9).Value = [double]"0.7238945645409351"
$ws.Cells.Item(2, 10).Value = [double]"0.7238945645409351"
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(2, 13).Value = [double]"0.003058333333333333"
$ws.Cells.Item(2, 14).Value = [double]"0.009175"
$ws.Cells.Item(2, 15).Value = [double]"0.0001134234803787887"
$ws.Cells.Item(2, 16).Value = [double]"0.0001134234803787887"
$ws.Cells.Item(2, 17).Value = [double]"0.112851408175"
$ws.Cells.Item(2, 18).Value = [double]"1.015662673575"
$ws.Cells.Item(2, 19).Value = [double]"8.210664093752054e-05"
$ws.Cells.Item(2, 20).Value = [double]"8.210664093752054e-05"

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efnb2"
$ws.Cells.Item(3, 3).Value = "Epha3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = [double]"36.899643"
$ws.Cells.Item(3, 8).Value = [double]"110.698929"
$ws.Cells.Item(3, 9).Value = [double]"0.7238945645409351"
$ws.Cells.Item(3, 10).Value = [double]"0.7238945645409351"
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = [double]"26.097779"
$ws.Cells.Item(3, 14).Value = [double]"78.29333700000001"
$ws.Cells.Item(3, 15).Value = [double]"0.9678804112271815"
$ws.Cells.Item(3, 16).Value = [double]"0.9678804112271815"
$ws.Cells.Item(3, 17).Value = [double]"962.998728192897"
$ws.Cells.Item(3, 18).Value = [double]"8666.988553736073"
$ws.Cells.Item(3, 19).Value = [double]"0.7006433688130018"
$ws.Cells.Item(3, 20).Value = [double]"0.7006433688130018"

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efnb2"
$ws.Cells.Item(4, 3).Value = "Epha3"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = [double]"36.899643"
$ws.Cells.Item(4, 8).Value = [double]"110.698929"
$ws.Cells.Item(4, 9).Value = [double]"0.7238945645409351"
$ws.Cells.Item(4, 10).Value = [double]"0.7238945645409351"
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(4, 13).Value = [double]"0.8630093333333333"
$ws.Cells.Item(4, 14).Value = [double]"2.589028"
$ws.Cells.Item(4, 15).Value = [double]"0.03200616529243972"
$ws.Cells.Item(4, 16).Value = [double]"0.03200616529243972"
$ws.Cells.Item(4, 17).Value = [double]"31.844736305668"
$ws.Cells.Item(4, 18).Value = [double]"286.602626751012"
$ws.Cells.Item(4, 19).Value = [double]"0.02316908908699584"
$ws.Cells.Item(4, 20).Value = [double]"0.02316908908699584"

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Efnb2"
$ws.Cells.Item(5, 3).Value = "Epha3"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = [double]"3.374819"
$ws.Cells.Item(5, 8).Value = [double]"10.124457"
$ws.Cells.Item(5, 9).Value = [double]"0.0662069584361419"
$ws.Cells.Item(5, 10).Value = [double]"0.0662069584361419"
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(5, 13).Value = [double]"0.003058333333333333"
$ws.Cells.Item(5, 14).Value = [double]"0.009175"
$ws.Cells.Item(5, 15).Value = [double]"0.0001134234803787887"
$ws.Cells.Item(5, 16).Value = [double]"0.0001134234803787887"
$ws.Cells.Item(5, 17).Value = [double]"0.01032132144166667"
$ws.Cells.Item(5, 18).Value = [double]"0.092891892975"
$ws.Cells.Item(5, 19).Value = [double]"7.509423651121018e-06"
$ws.Cells.Item(5, 20).Value = [double]"7.509423651121018e-06"

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efnb2"
$ws.Cells.Item(6, 3).Value = "Epha3"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = [double]"3.374819"
$ws.Cells.Item(6, 8).Value = [double]"10.124457"
$ws.Cells.Item(6, 9).Value = [double]"0.0662069584361419"
$ws.Cells.Item(6, 10).Value = [double]"0.0662069584361419"
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = [double]"26.097779"
$ws.Cells.Item(6, 14).Value = [double]"78.29333700000001"
$ws.Cells.Item(6, 15).Value = [double]"0.9678804112271815"
$ws.Cells.Item(6, 16).Value = [double]"0.9678804112271815"
$ws.Cells.Item(6, 17).Value = [double]"88.075280427001"
$ws.Cells.Item(6, 18).Value = [double]"792.677523843009"
$ws.Cells.Item(6, 19).Value = [double]"0.06408041815727393"
$ws.Cells.Item(6, 20).Value = [double]"0.06408041815727393"

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efnb2"
$ws.Cells.Item(7, 3).Value = "Epha3"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = [double]"3.374819"
$ws.Cells.Item(7, 8).Value = [double]"10.124457"
$ws.Cells.Item(7, 9).Value = [double]"0.0662069584361419"
$ws.Cells.Item(7, 10).Value = [double]"0.0662069584361419"
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(7, 13).Value = [double]"0.8630093333333333"
$ws.Cells.Item(7, 14).Value = [double]"2.589028"
$ws.Cells.Item(7, 15).Value = [double]"0.03200616529243972"
$ws.Cells.Item(7, 16).Value = [double]"0.03200616529243972"
$ws.Cells.Item(7, 17).Value = [double]"2.912500295310667"
$ws.Cells.Item(7, 18).Value = [double]"26.212502657796"
$ws.Cells.Item(7, 19).Value = [double]"0.002119030855216844"
$ws.Cells.Item(7, 20).Value = [double]"0.002119030855216844"

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Efnb2"
$ws.Cells.Item(8, 3).Value = "Epha3"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = [double]"10.699319"
$ws.Cells.Item(8, 8).Value = [double]"32.097957"
$ws.Cells.Item(8, 9).Value = [double]"0.2098984770229228"
$ws.Cells.Item(8, 10).Value = [double]"0.2098984770229228"
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(8, 13).Value = [double]"0.003058333333333333"
$ws.Cells.Item(8, 14).Value = [double]"0.009175"
$ws.Cells.Item(8, 15).Value = [double]"0.0001134234803787887"
$ws.Cells.Item(8, 16).Value = [double]"0.0001134234803787887"
$ws.Cells.Item(8, 17).Value = [double]"0.03272208394166667"
$ws.Cells.Item(8, 18).Value = [double]"0.294498755475"
$ws.Cells.Item(8, 19).Value = [double]"2.380741579014712e-05"
$ws.Cells.Item(8, 20).Value = [double]"2.380741579014712e-05"

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Efnb2"
$ws.Cells.Item(9, 3).Value = "Epha3"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = [double]"10.699319"
$ws.Cells.Item(9, 8).Value = [double]"32.097957"
$ws.Cells.Item(9, 9).Value = [double]"0.2098984770229228"
$ws.Cells.Item(9, 10).Value = [double]"0.2098984770229228"
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = [double]"26.097779"
$ws.Cells.Item(9, 14).Value = [double]"78.29333700000001"
$ws.Cells.Item(9, 15).Value = [double]"0.9678804112271815"
$ws.Cells.Item(9, 16).Value = [double]"0.9678804112271815"
$ws.Cells.Item(9, 17).Value = [double]"279.228462712501"
$ws.Cells.Item(9, 18).Value = [double]"2513.056164412509"
$ws.Cells.Item(9, 19).Value = [double]"0.2031566242569056"
$ws.Cells.Item(9, 20).Value = [double]"0.2031566242569056"

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Efnb2"
$ws.Cells.Item(10, 3).Value = "Epha3"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = [double]"10.699319"
$ws.Cells.Item(10, 8).Value = [double]"32.097957"
$ws.Cells.Item(10, 9).Value = [double]"0.2098984770229228"
$ws.Cells.Item(10, 10).Value = [double]"0.2098984770229228"
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(10, 13).Value = [double]"0.8630093333333333"
$ws.Cells.Item(10, 14).Value = [double]"2.589028"
$ws.Cells.Item(10, 15).Value = [double]"0.03200616529243972"
$ws.Cells.Item(10, 16).Value = [double]"0.03200616529243972"
$ws.Cells.Item(10, 17).Value = [double]"9.233612157310667"
$ws.Cells.Item(10, 18).Value = [double]"83.10250941579599"
$ws.Cells.Item(10, 19).Value = [double]"0.006718045350227029"
$ws.Cells.Item(10, 20).Value = [double]"0.006718045350227029"
